# The sheet contains a flat data table (header row 1, data rows 2-310).
# This edit inserts two brand-new records right before the current row 230,
# pushing the former rows 230-310 down to 232-312 (dimension grows from
# A1:R310 to A1:R312). Everything else in the sheet is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the affected block; Excel shifts every
# row at/after 230 down by two, carrying formatting (incl. the date number
# format on column D) along with it.
$ws.Rows("230:231").Insert()

# --- New row 230 -----------------------------------------------------------
$ws.Cells.Item(230, 1).Value  = 8
$ws.Cells.Item(230, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(230, 3).Value  = "Coquimbo"
$ws.Cells.Item(230, 4).Value  = 44795
$ws.Cells.Item(230, 5).Value  = 4
$ws.Cells.Item(230, 6).Value  = 100112003
$ws.Cells.Item(230, 7).Value  = "Ajo"
$ws.Cells.Item(230, 8).Value  = "Chino"
$ws.Cells.Item(230, 9).Value  = "Primera"
$ws.Cells.Item(230, 10).Value = 440
$ws.Cells.Item(230, 11).Value = 24500
$ws.Cells.Item(230, 12).Value = 25000
$ws.Cells.Item(230, 13).Value = 24750
$ws.Cells.Item(230, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(230, 15).Value = "China"
$ws.Cells.Item(230, 16).Value = 2475
$ws.Cells.Item(230, 17).Value = 10
$ws.Cells.Item(230, 18).Value = "Hortaliza"

# --- New row 231 -----------------------------------------------------------
$ws.Cells.Item(231, 1).Value  = 8
$ws.Cells.Item(231, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(231, 3).Value  = "Coquimbo"
$ws.Cells.Item(231, 4).Value  = 44795
$ws.Cells.Item(231, 5).Value  = 4
$ws.Cells.Item(231, 6).Value  = 100112003
$ws.Cells.Item(231, 7).Value  = "Ajo"
$ws.Cells.Item(231, 8).Value  = "Chino"
$ws.Cells.Item(231, 9).Value  = "Primera"
$ws.Cells.Item(231, 10).Value = 460
$ws.Cells.Item(231, 11).Value = 26500
$ws.Cells.Item(231, 12).Value = 27000
$ws.Cells.Item(231, 13).Value = 26750
$ws.Cells.Item(231, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(231, 15).Value = "China"
$ws.Cells.Item(231, 16).Value = 2675
$ws.Cells.Item(231, 17).Value = 10
$ws.Cells.Item(231, 18).Value = "Hortaliza"
